$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C3").Value = 157902
$ws.Range("C4").Value = 148959
$ws.Range("C5").Value = 8944
$ws.Range("C8").Value = 63.9
